$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.161.92"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = "'1.901.36"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.37%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'306.51"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  +1.51%  '
$ws.Range('D8').Value = "'0.3771"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.37%  '
$ws.Range('D9').Value = "'0.07246"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').Value = "'21.18"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.43%  '
$ws.Range('D11').Value = "'0.8986"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = "'0.08372"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +10.66%  '
$ws.Range('D13').Value = "'1.911.21"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('D14').Value = "'94.69"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = "'5.263"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').Value = "'0.000008584"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = "'14.49"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('D20').Value = "'27.204.32"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('D21').Value = "'5.055"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.54%  '
$ws.Range('D22').Value = "'2.150.39"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.14%  '
$ws.Range('D23').Value = "'10.59"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.90%  '
$ws.Range('D24').Value = "'6.419"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').Value = "'2.283"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.22%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'146.62"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.46%  '
$ws.Range('E27').Value = '  -1.44%  '
$ws.Range('D28').Value = "'18.10"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('D29').Value = "'114.59"
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Value = "'4.919"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('D31').Value = "'4.780"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.72%  '
$ws.Range('D32').Value = "'0.09225"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('D33').Value = "'0.8193"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.71%  '
$ws.Range('D34').Value = "'0.05048"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('E35').Value = '  +5.52%  '
$ws.Range('D36').Value = "'2.953"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.36%  '
$ws.Range('D37').Value = "'3.364"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.33%  '
$ws.Range('D38').Value = "'2.569"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.70%  '
$ws.Range('D39').Value = "'0.5689"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.06%  '
$ws.Range('D40').Value = "'0.01971"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.94%  '
$ws.Range('D41').Value = "'1.072"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').Value = "'6.655"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.36%  '
$ws.Range('D43').Value = "'8.933"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.94%  '
$ws.Range('D44').Value = "'118.28"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.87%  '
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').Value = "'0.4817"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').Value = "'10.18"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.97%  '
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D50').Value = "'37.43"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('D51').Value = "'63.56"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.36%  '
